# Edit script generated to reproduce the target diff for Variable_names_overview.xlsx
# (commit: "Versuch vom Proportion test")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old B29 cell (it held "table.ctrl1.SW.norm"); that text is
# reused/retitled below as the new D41 entry. ---
$ws.Cells.Item(29, 2).ClearContents()

# --- New block of rows 41-56 describing "Inhalt von ctrl_mean_filtert_final" ---

# Row 41 (header row of the new block)
$ws.Cells.Item(41, 4).Value = "Inhalt von ctrl_mean_filtert_final"
$ws.Cells.Item(41, 5).Value = "Bedeutung"

# Column D for rows 42-56, in row order
$ws.Cells.Item(42, 4).Value = "pb_fit"
$ws.Cells.Item(43, 4).Value = "sum_area"
$ws.Cells.Item(44, 4).Value = "fit_area"
$ws.Cells.Item(45, 4).Value = "fit_c_fxn"
$ws.Cells.Item(46, 4).Value = "fit_mean_fxn"
$ws.Cells.Item(47, 4).Value = "fit_param"
$ws.Cells.Item(48, 4).Value = "fitted"
$ws.Cells.Item(49, 4).Value = "fit_res"
$ws.Cells.Item(50, 4).Value = "fit_sigma"
$ws.Cells.Item(51, 4).Value = "fit_mean_fxn"
$ws.Cells.Item(52, 4).Value = "fit_c_fxn"
$ws.Cells.Item(53, 4).Value = "nb_max"
$ws.Cells.Item(54, 4).Value = "ctrl_max"
$ws.Cells.Item(55, 4).Value = "peaks"
$ws.Cells.Item(56, 4).Value = "maxima"

# Column E for the remaining rows, written in the same order the author
# originally typed them (preserves shared-string table ordering)
$ws.Cells.Item(53, 5).Value = "Anzahl lokale Maxima"
$ws.Cells.Item(55, 5).Value = "gibt es zusätzlich zum Hochpunkt ein Platteau"
$ws.Cells.Item(56, 5).Value = "Koordinate lokaler Maxima"
$ws.Cells.Item(54, 5).Value = "Koordinate lokaler Maxima und Platteaus"
$ws.Cells.Item(52, 5).Value = "Amplitude der Punkte"
$ws.Cells.Item(51, 5).Value = "Mittelwert der Punkte"
$ws.Cells.Item(50, 5).Value = "Standardabweichung der Punkte"
$ws.Cells.Item(49, 5).Value = "reseduals sum of squares"
$ws.Cells.Item(48, 5).Value = "Checkt, ob eine gefittete Kurve gefunden wurde"
$ws.Cells.Item(47, 5).Value = "Die Amplitude, Mittelwert und Standard Abweichung werden den Hochpunkten zugeordnet"
$ws.Cells.Item(46, 5).Value = "Mittelwert der Punkte nur Ränder werden mit einebzogen"
$ws.Cells.Item(45, 5).Value = "Amplitude der Punkte nur Ränder werden mit einebzogen"
$ws.Cells.Item(44, 5).Value = "Area unter jedem Maximum"
$ws.Cells.Item(43, 5).Value = " Summe der Areas"
$ws.Cells.Item(42, 5).Value = "Checkt, ob Kurve fittet"

# --- Column width adjustments for D (new wider content) and the newly used E column ---
$ws.Columns.Item(4).ColumnWidth = 30.75
$ws.Columns.Item(5).ColumnWidth = 25.92

# --- Update the view/selection to match where the author ended up editing ---
$ws.Range("D54").Select()
